$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.724.55'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.455.69'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.55'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.22'
$ws.Range('E6').Value = '  +7.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.455.49'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.65'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.391'
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.045.85'
$ws.Range('E13').Value = '  +2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.02'
$ws.Range('E14').Value = '  +8.71%  '
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.465.34'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.851.75'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  +8.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.39'
$ws.Range('E20').Value = '  +3.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.59'
$ws.Range('E21').Value = '  +2.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.73'
$ws.Range('E22').Value = '  +3.35%  '
$ws.Range('E23').Value = '  +2.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.70'
$ws.Range('E24').Value = '  +3.83%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.598.46'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.182'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.70'
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  -11.36%  '
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '24.18'
$ws.Range('E36').Value = '  +3.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.481.31'
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('E38').Value = '  +2.74%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '167.15'
$ws.Range('E41').Value = '  +1.58%  '
$ws.Range('E42').Value = '  +3.53%  '
$ws.Range('E43').Value = '  +6.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.806'
$ws.Range('E44').Value = '  +3.84%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.52'
$ws.Range('E45').Value = '  +3.97%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.59'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  +1.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.579.88'
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('E51').Value = '  +2.20%  '
